# Scheduled market-price refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) per leve row across all job sheets, per upstream data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 287.75
$ws.Range("N5").Value = -730
$ws.Range("L5").Value = 500
$ws.Range("I5").Value = 257.42856
$ws.Range("M5").Value = -142.42856
$ws.Range("K5").Value = 257.42856
$ws.Range("J5").Value = 500
$ws.Range("N62").Value = -5248
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("H62").Value = 4000
$ws.Range("N65").Value = -26240
$ws.Range("J65").Value = 4000
$ws.Range("H65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("M107").Value = 728.4286
$ws.Range("K107").Value = 1191.5714
$ws.Range("I107").Value = 1191.5714
$ws.Range("H107").Value = 954.65515
$ws.Range("N107").Value = -4172.75
$ws.Range("L107").Value = 332.75
$ws.Range("J107").Value = 332.75
$ws.Range("M132").Value = -2606.831
$ws.Range("L132").Value = 8299.5
$ws.Range("K132").Value = 5136.831
$ws.Range("N132").Value = -13359.5
$ws.Range("J132").Value = 2766.5
$ws.Range("H132").Value = 1801.3662
$ws.Range("I132").Value = 1712.277
$ws.Range("J137").Value = 2140
$ws.Range("N137").Value = -11520
$ws.Range("H137").Value = 1928.8422
$ws.Range("L137").Value = 6420
$ws.Range("H138").Value = 2440.0134
$ws.Range("M138").Value = -8173.699999999999
$ws.Range("K138").Value = 13313.7
$ws.Range("J138").Value = 2132.6462
$ws.Range("I138").Value = 4437.9
$ws.Range("N138").Value = -16677.9386
$ws.Range("L138").Value = 6397.9386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J137").Value = 37890
$ws.Range("N137").Value = -48090
$ws.Range("H137").Value = 41712
$ws.Range("L137").Value = 37890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4857.2666
$ws.Range("K134").Value = 12549.6
$ws.Range("I134").Value = 4183.2
$ws.Range("M134").Value = -10014.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M6").Value = -1667553.6
$ws.Range("L6").Value = 5500
$ws.Range("J6").Value = 5500
$ws.Range("I6").Value = 1667666.6
$ws.Range("K6").Value = 1667666.6
$ws.Range("H6").Value = 1002800
$ws.Range("N6").Value = -5726
$ws.Range("M7").Value = -157.25
$ws.Range("K7").Value = 270.25
$ws.Range("N7").Value = -342.833336
$ws.Range("I7").Value = 270.25
$ws.Range("H7").Value = 178.2
$ws.Range("J7").Value = 116.833336
$ws.Range("L7").Value = 116.833336
$ws.Range("M31").Value = -788.9749999999999
$ws.Range("I31").Value = 1083.975
$ws.Range("K31").Value = 1083.975
$ws.Range("N31").Value = -6142.3423
$ws.Range("L31").Value = 5552.3423
$ws.Range("H31").Value = 3260.8718
$ws.Range("J31").Value = 5552.3423
$ws.Range("I34").Value = 1083.975
$ws.Range("N34").Value = -5956.3423
$ws.Range("J34").Value = 5552.3423
$ws.Range("K34").Value = 1083.975
$ws.Range("L34").Value = 5552.3423
$ws.Range("M34").Value = -881.9749999999999
$ws.Range("H34").Value = 3260.8718
$ws.Range("L41").Value = 19999
$ws.Range("H41").Value = 13666
$ws.Range("M41").Value = -572
$ws.Range("J41").Value = 19999
$ws.Range("I41").Value = 1000
$ws.Range("N41").Value = -20855
$ws.Range("K41").Value = 1000
$ws.Range("M50").Value = -4458
$ws.Range("H50").Value = 15027
$ws.Range("K50").Value = 5083
$ws.Range("I50").Value = 5083
$ws.Range("I51").Value = 3090
$ws.Range("N51").Value = -19710.8
$ws.Range("L51").Value = 18238.8
$ws.Range("J51").Value = 18238.8
$ws.Range("H51").Value = 15714
$ws.Range("M51").Value = -2354
$ws.Range("K51").Value = 3090
$ws.Range("H59").Value = 27374.125
$ws.Range("L59").Value = 27374.125
$ws.Range("N59").Value = -29664.125
$ws.Range("J59").Value = 27374.125
$ws.Range("J60").Value = 10576
$ws.Range("K60").Value = 1093
$ws.Range("L60").Value = 10576
$ws.Range("M60").Value = -582
$ws.Range("I60").Value = 1093
$ws.Range("N60").Value = -11598
$ws.Range("H60").Value = 9390.625
$ws.Range("K61").Value = 3090
$ws.Range("H61").Value = 15714
$ws.Range("M61").Value = -2742
$ws.Range("L61").Value = 18238.8
$ws.Range("J61").Value = 18238.8
$ws.Range("N61").Value = -18934.8
$ws.Range("I61").Value = 3090
$ws.Range("L68").Value = 27555
$ws.Range("H68").Value = 27555
$ws.Range("N68").Value = -29053
$ws.Range("J68").Value = 27555
$ws.Range("L71").Value = 82665
$ws.Range("J71").Value = 27555
$ws.Range("N71").Value = -90153
$ws.Range("H71").Value = 27555
$ws.Range("J74").Value = 28438.3
$ws.Range("N74").Value = -30186.3
$ws.Range("H74").Value = 28438.3
$ws.Range("L74").Value = 28438.3
$ws.Range("J77").Value = 28438.3
$ws.Range("L77").Value = 85314.89999999999
$ws.Range("H77").Value = 28438.3
$ws.Range("N77").Value = -94050.89999999999
$ws.Range("H107").Value = 1002.7857
$ws.Range("N107").Value = -4945.9
$ws.Range("L107").Value = 1105.9
$ws.Range("J107").Value = 1105.9
$ws.Range("H134").Value = 2589.7078
$ws.Range("L134").Value = 4940.3079
$ws.Range("K134").Value = 8476.3272
$ws.Range("N134").Value = -10010.3079
$ws.Range("J134").Value = 1646.7693
$ws.Range("I134").Value = 2825.4424
$ws.Range("M134").Value = -5941.3272

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2218.2778
$ws.Range("N5").Value = -8066.8568
$ws.Range("L5").Value = 7842.8568
$ws.Range("I5").Value = 832.25
$ws.Range("M5").Value = -2384.75
$ws.Range("K5").Value = 2496.75
$ws.Range("J5").Value = 2614.2856
$ws.Range("H135").Value = 2218.2778
$ws.Range("I135").Value = 832.25
$ws.Range("L135").Value = 23528.5704
$ws.Range("K135").Value = 7490.25
$ws.Range("M135").Value = -4955.25
$ws.Range("N135").Value = -28598.5704
$ws.Range("J135").Value = 2614.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M113").Value = -77373.234
$ws.Range("H113").Value = 65178.875
$ws.Range("K113").Value = 79543.234
$ws.Range("I113").Value = 79543.234
$ws.Range("M132").Value = -8867.299999999999
$ws.Range("L132").Value = 16504.125
$ws.Range("K132").Value = 11397.3
$ws.Range("N132").Value = -21564.125
$ws.Range("J132").Value = 5501.375
$ws.Range("H132").Value = 4555.6665
$ws.Range("I132").Value = 3799.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M14").Value = -34828
$ws.Range("L14").Value = 54753.75
$ws.Range("I14").Value = 35000
$ws.Range("N14").Value = -55097.75
$ws.Range("H14").Value = 50803
$ws.Range("J14").Value = 54753.75
$ws.Range("K14").Value = 35000
$ws.Range("L122").Value = 12413.0772
$ws.Range("J122").Value = 4137.6924
$ws.Range("N122").Value = -17313.0772
$ws.Range("H122").Value = 4043.889
$ws.Range("I122").Value = 3800
$ws.Range("M122").Value = -8950
$ws.Range("K122").Value = 11400
$ws.Range("K136").Value = 6111.2607
$ws.Range("L136").Value = 100004886
$ws.Range("N136").Value = -100009986
$ws.Range("H136").Value = 5954345
$ws.Range("I136").Value = 2037.0869
$ws.Range("M136").Value = -3561.2607
$ws.Range("J136").Value = 33334962

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M9").Value = $null
$ws.Range("K9").Value = 0
$ws.Range("J9").Value = 67507
$ws.Range("H9").Value = 67507
$ws.Range("L9").Value = 67507
$ws.Range("N9").Value = -67787
$ws.Range("I9").Value = 0
$ws.Range("M14").Value = -33832
$ws.Range("L14").Value = 4891213
$ws.Range("I14").Value = 34000
$ws.Range("N14").Value = -4891549
$ws.Range("H14").Value = 3023054.2
$ws.Range("J14").Value = 4891213
$ws.Range("K14").Value = 34000
$ws.Range("K136").Value = 6755.286
$ws.Range("L136").Value = 11355.4614
$ws.Range("N136").Value = -16455.4614
$ws.Range("H136").Value = 2614.2
$ws.Range("I136").Value = 2251.762
$ws.Range("M136").Value = -4205.286
$ws.Range("J136").Value = 3785.1538
